$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2670
$ws.Range("D2").Value = 1357.5

$ws.Range("C3").Value = 2811
$ws.Range("D3").Value = 1492

$ws.Range("C4").Value = 2666
$ws.Range("D4").Value = 1347

$ws.Range("C5").Value = 2683
$ws.Range("D5").Value = 1382.5

$ws.Range("C6").Value = 2692
$ws.Range("D6").Value = 1390.5

$ws.Range("C7").Value = 2693
$ws.Range("D7").Value = 1383.5

$ws.Range("C8").Value = 2659
$ws.Range("D8").Value = 1337.5

$ws.Range("C9").Value = 2679
$ws.Range("D9").Value = 1371.5

$ws.Range("C10").Value = 2650
$ws.Range("D10").Value = 2650

$ws.Range("C11").Value = 2660
$ws.Range("D11").Value = 1337.5

$ws.Range("C12").Value = 2686.3
